# Insert a new "Supervisors" sheet before "Projects" and populate the
# supervisor-related data called for by the commit:
#   - a Supervisors sheet (Supervisor / Max_number_of_projects / Max_number_of_students)
#   - a Supervisor column added to the Projects sheet

$wb = $excel.ActiveWorkbook

$projectsRef = $wb.Worksheets.Item("Projects")

# --- New "Supervisors" sheet, inserted before "Projects" ---
$supervisors = $wb.Worksheets.Add($projectsRef)
$supervisors.Name = "Supervisors"

# Re-fetch "Projects" by name: the handle returned around `Add(Before)`
# can otherwise keep pointing at the newly inserted sheet.
$projects = $wb.Worksheets.Item("Projects")

$supervisors.Range("A1").Value = "Supervisor"
$supervisors.Range("B1").Value = "Max_number_of_projects"
$supervisors.Range("C1").Value = "Max_number_of_students"
$supervisors.Range("A2").Value = "Dr Smith"

$supervisors.Columns.Item(2).ColumnWidth = 22.54

# --- "Projects" sheet: rename header, add Supervisor column ---
$projects.Range("A1").Value = "Project"
$projects.Range("C1").Value = "Supervisor"
$projects.Range("C2").Value = "Dr Smith"
$projects.Range("C3").Value = "Dr Smith"
$projects.Range("C4").Value = "Dr Smith"
$projects.Range("C5").Value = "Dr Smith"
$projects.Range("C6").Value = "Dr Smith"

$projects.Columns.Item(2).ColumnWidth = 25.21

# "Projects" becomes the active sheet (diff: activeTab goes from 3 to 1,
# i.e. Supervisors=0, Projects=1)
$projects.Activate()
